# Auto-generated edit script: updates '想去人数' (F) and one '最低票价' (G) cell
# per the commit diff, across sheets 展览(1), 演出(2), 全部类型(4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 63
$ws.Range("F4").Value = 1336
$ws.Range("F6").Value = 620
$ws.Range("F7").Value = 894
$ws.Range("F8").Value = 1546
$ws.Range("F9").Value = 160
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = 1452
$ws.Range("F12").Value = 3079
$ws.Range("F13").Value = 612
$ws.Range("F14").Value = 1763
$ws.Range("F15").Value = 1798
$ws.Range("F16").Value = 847
$ws.Range("F17").Value = 273
$ws.Range("G18").Value = 168
$ws.Range("F19").Value = 1467
$ws.Range("F20").Value = 291
$ws.Range("F22").Value = 6
$ws.Range("F23").Value = 1206
$ws.Range("F24").Value = 400
$ws.Range("F25").Value = 450
$ws.Range("F26").Value = 104
$ws.Range("F27").Value = 4765
$ws.Range("F28").Value = 29
$ws.Range("F31").Value = 1643
$ws.Range("F33").Value = 132
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 55
$ws.Range("F3").Value = 29
$ws.Range("F6").Value = 54
$ws.Range("F7").Value = 67
$ws.Range("F14").Value = 2
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 63
$ws.Range("F4").Value = 55
$ws.Range("F5").Value = 29
$ws.Range("F9").Value = 54
$ws.Range("F10").Value = 67
$ws.Range("F12").Value = 1336
$ws.Range("F14").Value = 620
$ws.Range("F15").Value = 895
$ws.Range("F16").Value = 1546
$ws.Range("F17").Value = 160
$ws.Range("F18").Value = 160
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 1452
$ws.Range("F22").Value = 3079
$ws.Range("F23").Value = 612
$ws.Range("F24").Value = 1763
$ws.Range("F25").Value = 1798
$ws.Range("F26").Value = 847
$ws.Range("F27").Value = 273
$ws.Range("G28").Value = 168
$ws.Range("F29").Value = 1467
$ws.Range("F30").Value = 291
$ws.Range("F33").Value = 6
$ws.Range("F35").Value = 1206
$ws.Range("F36").Value = 400
$ws.Range("F37").Value = 450
$ws.Range("F38").Value = 104
$ws.Range("F39").Value = 4765
$ws.Range("F40").Value = 29
$ws.Range("F43").Value = 1643
$ws.Range("F47").Value = 132
$ws.Range("F48").Value = 2

